$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hx")
Write-Host $ws.Name
